# Added script for cash payment
#
# - test_suite: adds a new "MarkPaymentAsPaid" test case, renames
#   PayThruUnionBank -> PayThruUnionPay, and flips a few runmode flags.
# - BillingAddress: adds a "Reference Number" column.
# - Updates the active sheet / selections to match the authored state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("test_suite")
$ws2 = $wb.Worksheets.Item("OrderDetails")
$ws3 = $wb.Worksheets.Item("BillingAddress")

# --- BillingAddress sheet: new "Reference Number" column --------------
$ws3.Range("N1").Value = "Reference Number"

# --- test_suite sheet ---------------------------------------------------
# Rename PayThruUnionBank -> PayThruUnionPay (still row 10 at this point,
# before the new row is inserted below).
$ws1.Range("A10").Value = "PayThruUnionPay"

# Insert a new row for the "MarkPaymentAsPaid" test case right after
# PayThruCash (was row 2), pushing everything below it down by one.
$ws1.Rows.Item(3).Insert()
$ws1.Range("A3").Value = "MarkPaymentAsPaid"
$ws1.Range("B3").Value = "N"

# Update the runmode flags that changed.
$ws1.Range("B2").Value = "N"    # PayThruCash: Y -> N
$ws1.Range("B4").Value = "N"    # PayThruGcash: Y -> N
$ws1.Range("B8").Value = "Y"    # PayThruAliPay: N -> Y

# --- View state: active sheet + selections ------------------------------
# Selection on test_suite moves to B7.
$ws1.Activate()
$ws1.Range("B7").Select()

# Selection on BillingAddress moves to G3 (and view scrolls right).
$ws3.Activate()
$ws3.Range("G3").Select()

# OrderDetails becomes the active sheet/tab, with selection at E17.
$ws2.Activate()
$ws2.Range("E17").Select()
